$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = -2
    4  = 1
    5  = -3
    6  = -3
    7  = 1
    8  = 3
    9  = 2
    10 = 1
    11 = -5
    12 = -1
    14 = -1
    16 = 1
    17 = 3
    18 = 1
    19 = 4
    20 = -5
    21 = 3
    22 = -8
    23 = -5
    25 = -1
    29 = 1
    30 = -5
    31 = -9
    33 = -2
    36 = 1
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 6).Value = $values[$row]
}
